# Swap the order of the "Recorded By" names in column G:
# "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$colG = $ws.Columns.Item(7)

$null = $colG.Replace(
    "dnasr281@gmail.com, System",  # What
    "System, dnasr281@gmail.com",  # Replacement
    1,                             # LookAt: xlWhole
    1,                             # SearchOrder: xlByRows
    $false,                        # MatchCase
    $false,                        # MatchByte
    $true                          # SearchFormat
)
